# ---------------------------------------------------------------------------
# Edit script for "Crystal Palace_stats.xlsx"
#
# 1) Rename stat sheets to human-friendly, spaced-out names.
# 2) Bump every player's "Age" column (format YY-DDD, years-days) forward
#    by one day on every per-player stats sheet.
# 3) Fix the "Playing Time" merged header on StandardStats / PlayingTime:
#    it was incorrectly merged starting at F1 (leaving an extra blank
#    "Unnamed: 4_level_0" column hidden) - it should start at G1, with F1
#    holding its own "Unnamed: 4_level_0" label.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename worksheets
# ---------------------------------------------------------------------------
$renames = @(
    @{ Old = "StandardStats";    New = "Standard Stats" },
    @{ Old = "ShootingStats";    New = "Shooting Stats" },
    @{ Old = "PassingStats";     New = "Passing Stats" },
    @{ Old = "PassTypes";        New = "Pass Types" },
    @{ Old = "GoalShotCreation"; New = "Goal & Shot Creation" },
    @{ Old = "DefensiveActions"; New = "Defensive Actions" },
    @{ Old = "PlayingTime";      New = "Playing Time" },
    @{ Old = "MiscStats";        New = "Miscellaneous Stats" }
)

foreach ($r in $renames) {
    $sheet = $wb.Worksheets.Item($r.Old)
    $sheet.Name = $r.New
}

# ---------------------------------------------------------------------------
# 2) Age column (column E, rows 4-30 / 4-40) bumped by one day
# ---------------------------------------------------------------------------

$ageUpdatesFull = @(
    @{ Row = 4; New = "28-009" },
    @{ Row = 5; New = "25-201" },
    @{ Row = 6; New = "24-251" },
    @{ Row = 7; New = "28-299" },
    @{ Row = 8; New = "24-349" },
    @{ Row = 9; New = "27-266" },
    @{ Row = 10; New = "26-265" },
    @{ Row = 11; New = "27-024" },
    @{ Row = 12; New = "29-338" },
    @{ Row = 13; New = "30-147" },
    @{ Row = 14; New = "24-358" },
    @{ Row = 15; New = "25-259" },
    @{ Row = 16; New = "28-228" },
    @{ Row = 17; New = "20-292" },
    @{ Row = 18; New = "25-295" },
    @{ Row = 19; New = "25-072" },
    @{ Row = 20; New = "33-350" },
    @{ Row = 21; New = "21-161" },
    @{ Row = 22; New = "27-064" },
    @{ Row = 23; New = "28-294" },
    @{ Row = 24; New = "21-277" },
    @{ Row = 25; New = "32-088" },
    @{ Row = 26; New = "28-090" },
    @{ Row = 27; New = "33-191" },
    @{ Row = 28; New = "18-249" },
    @{ Row = 29; New = "19-312" },
    @{ Row = 30; New = "19-174" },
    @{ Row = 31; New = "22-357" },
    @{ Row = 32; New = "20-354" },
    @{ Row = 33; New = "29-182" },
    @{ Row = 34; New = "31-361" },
    @{ Row = 35; New = "19-166" },
    @{ Row = 36; New = "31-039" },
    @{ Row = 37; New = "20-165" },
    @{ Row = 38; New = "30-270" },
    @{ Row = 39; New = "20-054" },
    @{ Row = 40; New = "35-143" }
)

$ageUpdatesPartial = @(
    @{ Row = 4; New = "28-009" },
    @{ Row = 5; New = "25-201" },
    @{ Row = 6; New = "24-251" },
    @{ Row = 7; New = "28-299" },
    @{ Row = 8; New = "24-349" },
    @{ Row = 9; New = "27-266" },
    @{ Row = 10; New = "26-265" },
    @{ Row = 11; New = "27-024" },
    @{ Row = 12; New = "29-338" },
    @{ Row = 13; New = "30-147" },
    @{ Row = 14; New = "24-358" },
    @{ Row = 15; New = "25-259" },
    @{ Row = 16; New = "28-228" },
    @{ Row = 17; New = "20-292" },
    @{ Row = 18; New = "25-295" },
    @{ Row = 19; New = "25-072" },
    @{ Row = 20; New = "33-350" },
    @{ Row = 21; New = "21-161" },
    @{ Row = 22; New = "27-064" },
    @{ Row = 23; New = "28-294" },
    @{ Row = 24; New = "21-277" },
    @{ Row = 25; New = "32-088" },
    @{ Row = 26; New = "28-090" },
    @{ Row = 27; New = "33-191" },
    @{ Row = 28; New = "18-249" },
    @{ Row = 29; New = "19-312" },
    @{ Row = 30; New = "19-174" }
)

# Sheets that list 37 players (rows 4-40): Standard Stats, Playing Time
$fullSheets = @("Standard Stats", "Playing Time")

# Sheets that list 27 players (rows 4-30)
$partialSheets = @(
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Miscellaneous Stats"
)

foreach ($sheetName in $fullSheets) {
    $sheet = $wb.Worksheets.Item($sheetName)
    foreach ($u in $ageUpdatesFull) {
        $sheet.Cells.Item($u.Row, 5).Value = $u.New
    }
}

foreach ($sheetName in $partialSheets) {
    $sheet = $wb.Worksheets.Item($sheetName)
    foreach ($u in $ageUpdatesPartial) {
        $sheet.Cells.Item($u.Row, 5).Value = $u.New
    }
}

# ---------------------------------------------------------------------------
# 3) Fix merged "Playing Time" header (F1:I1 -> G1:I1) on the two sheets
#    that carry this banner, giving F1 its own "Unnamed: 4_level_0" label.
# ---------------------------------------------------------------------------
$headerFixSheets = @("Standard Stats", "Playing Time")

foreach ($sheetName in $headerFixSheets) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range("F1:I1").UnMerge()
    $sheet.Range("F1").Value = "Unnamed: 4_level_0"
    $sheet.Range("G1").Value = "Playing Time"
    $sheet.Range("G1:I1").Merge()
}

Write-Host "Edit complete."
